$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the import text in B2: remove the stray line break and fix "Com." -> "com."
$ws.Range("B2").Value = "com.blackknight.demo.models.MortgageRequest,`ncom.blackknight.demo.models.Address,com.blackknight.demo.models.Loan"

# 2. Turn off wrap text for D5:G5 (this removes the redundant style and
#    causes the cells to pick up the existing wrap-text style slot).
$ws.Range("D5:G5").WrapText = $false

# 3. Adjust column widths
$ws.Columns.Item(2).ColumnWidth = 61.92
$ws.Columns.Item(3).ColumnWidth = 22.7

# 4. Adjust row height for row 2
$ws.Rows.Item(2).RowHeight = 22.8

# 5. Change the active cell / view location
$ws.Range("A1").Select()
$ws.Range("B4").Select()
